# edit.ps1 -- applies the k_3.xlsx diff:
#  * add sheet 'o_20' after 'o_10'
#  * add sheet 'o_20_jumbled' after 'o_20'
#  * add column E 'evaluator_partial_correctness' to 'o_10' (with header style copied from D1)
#  * populate each sheet's header row + data row 2
#  * re-word each prompt's inline questions to append
#    'Return the sequence of nodes in response.'
#  * change evaluator_response on o_10 from 'Correct' to 'invalid input'
#  * keep 'o_10' as the active/selected sheet at the end

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "o_10"

# ===== Text blocks (single-quoted here-strings -> fully literal, no interpolation) =====
$promptO10 = @'
 Given is the adjacency matrix for a unweighted undirected graph containing 10 nodes labelled A to J. The value corresponding to each row M and column N represents whether there is a connection between the two nodes, where 0 means no connection.   
Consider some examples
Example 1: what is the shortest path from node A to node K? Return the sequence of nodes in response.
   A B C D E F G H I J K
 A 0 1 0 0 0 0 0 0 0 0 0
 B 1 0 1 0 0 0 0 0 0 0 0
 C 0 1 0 1 0 0 0 0 0 0 0
 D 0 0 1 0 1 0 0 0 0 0 0
 E 0 0 0 1 0 1 0 0 0 0 0
 F 0 0 0 0 1 0 1 0 0 0 0
 G 0 0 0 0 0 1 0 1 0 0 0
 H 0 0 0 0 0 0 1 0 1 0 0
 I 0 0 0 0 0 0 0 1 0 1 0
 J 0 0 0 0 0 0 0 0 1 0 1
 K 0 0 0 0 0 0 0 0 0 1 0
Solution: A -> B -> C -> D -> E -> F -> G -> H -> I -> J -> K
Example 2: what is the shortest path from node A to node F? Return the sequence of nodes in response.
   A B C D E F
 A 0 1 0 0 0 0
 B 1 0 1 0 0 0
 C 0 1 0 1 0 0
 D 0 0 1 0 1 0
 E 0 0 0 1 0 1
 F 0 0 0 0 1 0
Solution: A -> B -> C -> D -> E -> F
Example 3: what is the shortest path from node A to node I? Return the sequence of nodes in response.
   A B C D E F G H I
 A 0 1 0 0 0 0 0 0 0
 B 1 0 1 0 0 0 0 0 0
 C 0 1 0 1 0 0 0 0 0
 D 0 0 1 0 1 0 0 0 0
 E 0 0 0 1 0 1 0 0 0
 F 0 0 0 0 1 0 1 0 0
 G 0 0 0 0 0 1 0 1 0
 H 0 0 0 0 0 0 1 0 1
 I 0 0 0 0 0 0 0 1 0
Solution: A -> B -> C -> D -> E -> F -> G -> H -> I
 Given these examples, answer the following quesiton.
what is the shortest path from node A to node J? Return the sequence of nodes in response.
   A B C D E F G H I J
 A 0 1 0 0 0 0 0 0 0 0
 B 1 0 1 0 0 0 0 0 0 0
 C 0 1 0 1 0 0 0 0 0 0
 D 0 0 1 0 1 0 0 0 0 0
 E 0 0 0 1 0 1 0 0 0 0
 F 0 0 0 0 1 0 1 0 0 0
 G 0 0 0 0 0 1 0 1 0 0
 H 0 0 0 0 0 0 1 0 1 0
 I 0 0 0 0 0 0 0 1 0 1
 J 0 0 0 0 0 0 0 0 1 0
    
'@
$solutionO10 = @'
A -> B -> C -> D -> E -> F -> G -> H -> I -> J
'@
$llmResponseO10 = @'
The shortest path from node A to node J is: A -> B -> C -> D -> E -> F -> G -> H -> I -> J
'@
$evaluatorResponse = @'
invalid input
'@
$evalPartialO10 = @'
10/10
'@
$promptO20 = @'
 Given is the adjacency matrix for a unweighted undirected graph containing 20 nodes labelled A to T. The value corresponding to each row M and column N represents whether there is a connection between the two nodes, where 0 means no connection.   
Consider some examples
Example 1: what is the shortest path from node A to node X? Return the sequence of nodes in response.
   A B C D E F G H I J K L M N O P Q R S T U V W X
 A 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 G 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 H 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 I 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0
 K 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0
 M 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0
 R 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0
 S 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0
 T 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0
 U 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0
 V 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0
 W 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1
 X 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0
Solution: A -> B -> C -> D -> E -> F -> G -> H -> I -> J -> K -> L -> M -> N -> O -> P -> Q -> R -> S -> T -> U -> V -> W -> X
Example 2: what is the shortest path from node A to node S? Return the sequence of nodes in response.
   A B C D E F G H I J K L M N O P Q R S
 A 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0
 G 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0
 H 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0
 I 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0
 K 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0
 M 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0
 R 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1
 S 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0
Solution: A -> B -> C -> D -> E -> F -> G -> H -> I -> J -> K -> L -> M -> N -> O -> P -> Q -> R -> S
Example 3: what is the shortest path from node A to node Q? Return the sequence of nodes in response.
   A B C D E F G H I J K L M N O P Q
 A 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0
 F 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0
 G 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0
 H 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0
 I 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0
 J 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0
 K 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0
 L 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0
 M 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0
 N 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0
 O 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0
 P 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1
 Q 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0
Solution: A -> B -> C -> D -> E -> F -> G -> H -> I -> J -> K -> L -> M -> N -> O -> P -> Q
 Given these examples, answer the following quesiton.
what is the shortest path from node A to node T? Return the sequence of nodes in response.
   A B C D E F G H I J K L M N O P Q R S T
 A 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0
 G 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0
 H 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0
 I 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0
 K 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0
 M 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0
 R 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0
 S 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1
 T 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0
    
'@
$solutionO20 = @'
A -> B -> C -> D -> E -> F -> G -> H -> I -> J -> K -> L -> M -> N -> O -> P -> Q -> R -> S -> T
'@
$llmResponseO20 = @'
The shortest path from node A to node T is: A -> B -> C -> D -> E -> F -> G -> H -> I -> J -> K -> L -> M -> N -> O -> P -> Q -> R -> S -> T
'@
$evalPartialO20 = @'
20/20
'@
$llmResponseO20Jumbled = @'
The shortest path from node A to node T is: A -> B -> C -> D -> E -> F -> G -> H -> I -> J -> K -> L -> M -> N -> O -> P -> Q -> R -> S -> T.
'@

# ===== 1. Update o_10 (sheet1): rewrite prompt/solution/llm_response/evaluator_response, =====
# =====    and add the new "evaluator_partial_correctness" column E.                       =====
$ws1.Range("A2").Value = $promptO10
$ws1.Range("B2").Value = $solutionO10
$ws1.Range("C2").Value = $llmResponseO10
$ws1.Range("D2").Value = $evaluatorResponse

# Give E1 the same header formatting (bold, border, centered) as the existing header cells,
# by copying the format from D1 before writing the new header text/value.
$ws1.Range("D1").Copy() | Out-Null
$ws1.Range("E1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws1.Range("E1").Value = "evaluator_partial_correctness"
$ws1.Range("E2").Value = $evalPartialO10
$ws1.Rows.Item(2).EntireRow.AutoFit() | Out-Null

# ===== 2. Add sheet "o_20" right after "o_10" =====
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "o_20"

$ws1.Range("A1:E1").Copy() | Out-Null
$ws2.Range("A1:E1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws2.Range("A1").Value = "prompt"
$ws2.Range("B1").Value = "solution"
$ws2.Range("C1").Value = "llm_response"
$ws2.Range("D1").Value = "evaluator_response"
$ws2.Range("E1").Value = "evaluator_partial_correctness"

$ws2.Range("A2").Value = $promptO20
$ws2.Range("B2").Value = $solutionO20
$ws2.Range("C2").Value = $llmResponseO20
$ws2.Range("D2").Value = $evaluatorResponse
$ws2.Range("E2").Value = $evalPartialO20
$ws2.Rows.Item(2).EntireRow.AutoFit() | Out-Null

# ===== 3. Add sheet "o_20_jumbled" right after "o_20" =====
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "o_20_jumbled"

$ws1.Range("A1:E1").Copy() | Out-Null
$ws3.Range("A1:E1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws3.Range("A1").Value = "prompt"
$ws3.Range("B1").Value = "solution"
$ws3.Range("C1").Value = "llm_response"
$ws3.Range("D1").Value = "evaluator_response"
$ws3.Range("E1").Value = "evaluator_partial_correctness"

$ws3.Range("A2").Value = $promptO20
$ws3.Range("B2").Value = $solutionO20
$ws3.Range("C2").Value = $llmResponseO20Jumbled
$ws3.Range("D2").Value = $evaluatorResponse
$ws3.Range("E2").Value = $evalPartialO20
$ws3.Rows.Item(2).EntireRow.AutoFit() | Out-Null

# ===== 4. Restore "o_10" as the active/selected sheet (matches tabSelected="1" on sheet1 only) =====
$ws1.Activate()
$ws1.Range("A1").Select() | Out-Null

Write-Host "Sheets now: $($wb.Worksheets.Count) -> $([string]::Join(', ', @(1..$wb.Worksheets.Count | ForEach-Object { $wb.Worksheets.Item($_).Name })))"
